$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grievanceDetails")

# The grievance sample data was changed from a "Public Health and Sanitation /
# Mosquito menace" complaint to a "Street Lighting / Non Burning of Street
# Lights" complaint (Forward/Close Grievance test data).
$ws.Range("B2").Value = "Street Lighting"
$ws.Range("C2").Value = "Non Burning of Street Lights"
$ws.Range("D2").Value = "No street light past 3 days"

# Columns C & D are narrower/wider respectively to fit the new text.
$ws.Columns.Item(3).ColumnWidth = 25.5
$ws.Columns.Item(4).ColumnWidth = 23.3

# Move the active selection (the sheet is no longer scrolled to column E).
[void]$ws.Range("C6").Select()
